$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: A1 = "relation", B1 = "count"
$ws.Range("A1").Value = "relation"
$ws.Range("B1").Value = "count"

# Set column A width (target raw OOXML width 54.1640625; COM ColumnWidth
# quantizes to 1/6 steps, so 53.33 is the closest input that lands on the
# nearest achievable raw width of 54.166666666666664)
$ws.Columns.Item(1).ColumnWidth = 53.33

# Set selection to B1
$ws.Range("B1").Select()
